$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Jidlo" (food) column text, row by row (top to bottom) -
# same dish, new label text.
$ws.Range("C2").Value = "Food 111/ Traditional Hungarian goulash with pork meat and sauerkraut served with bread dumplings"
$ws.Range("C3").Value = "Food 222 / Old-czech style turkey breast with jasmine rice"
$ws.Range("C4").Value = "Foood 333 / Grilled minced meat with roasted potatoes and spicy salad of roasted peppers"
$ws.Range("C5").Value = "Food 444 / Tagliolini with beef tenderloin sprinkled with Grana Padano Cheese"
$ws.Range("C6").Value = "Foood 555  / Baked zander with vegetables in butter served with parsley potatoes"

# Update the "Polivka" (soup) column text, row by row (bottom to top) -
# same dish, new label text.
$ws.Range("F6").Value = "5555 polívka / Bank holiday. We do not serve daily menu."
$ws.Range("F5").Value = "4444 polívka / Minestrone soup with pasta"
$ws.Range("F4").Value = "33333 / Beef consommé with meat and noodles"
$ws.Range("F3").Value = "2222/ Lentil soup with sausages"
$ws.Range("F2").Value = "Mlsná 111 polévka 111 / Potato soup"

# Shift the week's date range forward by one week (tomorrow's / next week's orders).
$ws.Range("A2").Value = "7/28/2025"
$ws.Range("B2").Value = "8/1/2025"

# Move the active selection to L2.
$ws.Range("L2").Select()
